$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New employee rows to append below the existing data (rows 2-6 already exist)
$rows = @(
    @(306, "Peter",   "Parker",    40000, 300),
    @(307, "Mary",    "Jane",      50000, 306),
    @(308, "Tony",    "Stark",     70000, 300),
    @(310, "Steve",   "Rogers",    40000, 308),
    @(311, "Natasha", "Romanoff",  45000, 310),
    @(450, "Nick",    "Fury",      50000, 300),
    @(353, "Stephen", "Strange",   35000, 450),
    @(367, "Bruce",   "Banner",    40000, 450)
)

$r = 7
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# Match the column C autofit ("best fit") width shown in the diff as closely
# as this host's character/pixel-grid ColumnWidth quantization allows.
$ws.Columns.Item(3).ColumnWidth = 8.25

# Match the resulting selection from the diff
$ws.Range("D13").Select()
